$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60
$ws.Range("B60").Value = 0.070833333333333331
$ws.Range("B60").NumberFormat = "h:mm"
$ws.Range("C60").Value = 0.078472222222222221
$ws.Range("C60").NumberFormat = "h:mm"
$ws.Range("E60").Value = 11
$ws.Range("F60").Value = "Code"
$ws.Range("G60").Value = "Formatting the new order confirmation dialog pop-up box"

# Row 61
$ws.Range("B61").Value = 0.095138888888888884
$ws.Range("B61").NumberFormat = "h:mm"
$ws.Range("G61").Value = "Get total value amounts from all products in order and add to order"
$ws.Range("C61").Value = "2:53PM"
$ws.Range("E61").Value = 36
$ws.Range("F61").Value = "Code"
$ws.Range("H61").Value = "Orders now get added to DB. Now need to add prods too"

# Row 62
$ws.Range("B62").Value = "2:53PM"
$ws.Range("F62").Value = "Debugging"
$ws.Range("G62").Value = "Debugging sql inserts and queries for newOrder and products in new order"

# Column F width change (target stored width 25.28515625; this runtime quantizes
# column widths to steps of 1/6, so 24.45 is the closest input that lands on the
# nearest representable stored width of 25.333333333333332)
$ws.Columns("F").ColumnWidth = 24.45

# View changes
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G62").Select()
